$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The H column used a custom boolean-style number format ("TRUE";"TRUE";"FALSE")
# driving an =FALSE() formula. The author replaced that with the literal text
# "False" in every data row, reformatting the column as plain Text and leaving
# a few extra blank (but still Text-formatted) rows below the table.

# Re-format the whole H column range (existing rows + the new blank rows)
# as Text so the values stick as literal strings instead of booleans.
$ws.Range("H2:H40").NumberFormat = "@"

# Replace every =FALSE() cell with the literal text "False" (leading
# apostrophe forces text entry so Excel doesn't re-coerce "False" back into
# a boolean).
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 8).Value = "'False"
}

# Leave H38:H40 present but empty (just Text-formatted), matching the
# extra rows appended below the table.

# Update the view: scrolled down a bit with a new selection below the table.
$win = $excel.ActiveWindow
$ws.Range("G38:I41").Select()
$win.Zoom = 100
